$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing return/volatility rows (row3=return, row4=volatility) ---
$ws.Range("E3").Value = 0.06
$ws.Range("F3").Value = 0.08

$ws.Range("B4").Value = 0.08
$ws.Range("C4").Value = 0.12
$ws.Range("D4").Value = 0.04
$ws.Range("E4").Value = 0.04
$ws.Range("F4").Value = 0.12
$ws.Range("G4").Value = 0.2

# --- Add new "risk_free" row (row 5) ---
$ws.Range("A5").Value = "risk_free"
$ws.Range("B5").Value = 0.03
$ws.Range("C5").Value = 0.03
$ws.Range("D5").Value = 0.01
$ws.Range("E5").Value = 0.05
$ws.Range("F5").Value = 0.05
$ws.Range("G5").Value = 0.05
$ws.Range("B5:G5").Style = $ws.Range("B4").Style

# --- Add "sharpe_ratio" row (row 6) with formula ---
$ws.Range("A6").Value = "sharpe_ratio"
$ws.Range("B6:G6").Formula = "=(B3-B5)/B4"
$ws.Range("B6:G6").Style = $ws.Range("B3").Style
$ws.Range("B6:G6").NumberFormat = "0.0%"

# --- Remove old row 7 (no longer present) ---
$ws.Rows("7").Delete()

# --- Column width for column A ---
$ws.Columns("A").ColumnWidth = 10.69140625

# --- Selection ---
$ws.Range("A6").Select()
